# "Added New Mac-Address and Document Types"
# Append 5 new reg-center/machine/device rows (regcntr_id 10002 / machine_id 10032,
# device_ids 3000176-3000180) to the master-reg_center_machine_device table, and
# switch the workbook to manual calculation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
if (-not $ws) { $ws = $wb.Worksheets.Item(1) }

# Workbook now opens in manual calculation mode (calcPr@calcMode="manual")
$excel.Calculation = -4135

$startRow = 157
$regcntrId = 10002
$machineId = 10032
$deviceIds = 3000176, 3000177, 3000178, 3000179, 3000180

for ($i = 0; $i -lt $deviceIds.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $regcntrId      # regcntr_id
    $ws.Cells.Item($r, 2).Value = $machineId       # machine_id
    $ws.Cells.Item($r, 3).Value = $deviceIds[$i]   # device_id
    $ws.Cells.Item($r, 4).Value = "eng"            # lang_code
    $ws.Cells.Item($r, 5).Value = $true            # is_active
    $ws.Cells.Item($r, 6).Value = "superadmin"     # cr_by
    $ws.Cells.Item($r, 7).Value = "now()"          # cr_dtimes
    $ws.Cells.Item($r, 8).Value = "now()"          # eff_dtimes
}

# Leave the selection on the last-entered row's lang_code cell, as in the source edit.
$ws.Range("D157").Select() | Out-Null
